$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 111708029
$ws.Range("B2").Value = 90662
$ws.Range("D2").Value = 'LC'
$ws.Range("E2").Value = 4363
$ws.Range("F2").Value = 'Zontaggsvamp'
$ws.Range("G2").Value = 'Hydnellum concrescens'
$ws.Range("H2").Value = '(Pers.) Banker'
$ws.Range("Z2").Value = '14:21'
$ws.Range("AB2").Value = '14:21'

# Row 3
$ws.Range("A3").Value = 111706580
$ws.Range("B3").Value = 88032
$ws.Range("D3").Value = 'VU'
$ws.Range("E3").Value = 6276
$ws.Range("F3").Value = 'Goliatmusseron'
$ws.Range("G3").Value = 'Tricholoma matsutake'
$ws.Range("H3").Value = '(S.Ito & S.Imai) Singer'
$ws.Range("Z3").Value = '14:48'
$ws.Range("AB3").Value = '14:48'

# Row 4
$ws.Range("A4").Value = 111708162
$ws.Range("B4").Value = 90658
$ws.Range("D4").Value = 'NT'
$ws.Range("E4").Value = 4361
$ws.Range("F4").Value = 'Orange taggsvamp'
$ws.Range("G4").Value = 'Hydnellum aurantiacum'
$ws.Range("H4").Value = '(Batsch:Fr.) P.Karst.'
$ws.Range("Z4").Value = '14:12'
$ws.Range("AB4").Value = '14:12'

# Row 5
$ws.Range("A5").Value = 111708126
$ws.Range("B5").Value = 90666
$ws.Range("D5").Value = 'LC'
$ws.Range("E5").Value = 4364
$ws.Range("F5").Value = 'Dropptaggsvamp'
$ws.Range("G5").Value = 'Hydnellum ferrugineum'
$ws.Range("H5").Value = '(Fr.:Fr.) P. Karst.'
$ws.Range("Z5").Value = '14:14'
$ws.Range("AB5").Value = '14:14'

# Row 6
$ws.Range("A6").Value = 111704319
$ws.Range("B6").Value = 90710
$ws.Range("D6").Value = 'NT'
$ws.Range("E6").Value = 5449
$ws.Range("F6").Value = 'Svart taggsvamp'
$ws.Range("G6").Value = 'Phellodon niger'
$ws.Range("H6").Value = '(Fr.:Fr.) P.Karst.'
$ws.Range("Z6").Value = '15:11'
$ws.Range("AB6").Value = '15:11'

# Row 7
$ws.Range("A7").Value = 111708888
$ws.Range("B7").Value = 90678
$ws.Range("D7").Value = 'LC'
$ws.Range("E7").Value = 4366
$ws.Range("F7").Value = 'Skarp dropptaggsvamp'
$ws.Range("G7").Value = 'Hydnellum peckii'
$ws.Range("H7").Value = 'Banker'
$ws.Range("Z7").Value = '13:54'
$ws.Range("AB7").Value = '13:54'

# Row 8
$ws.Range("A8").Value = 111708099
$ws.Range("B8").Value = 90660
$ws.Range("D8").Value = 'NT'
$ws.Range("E8").Value = 4362
$ws.Range("F8").Value = 'Blå taggsvamp'
$ws.Range("G8").Value = 'Hydnellum caeruleum'
$ws.Range("H8").Value = '(Hornem.) P.Karst.'
$ws.Range("Z8").Value = '14:16'
$ws.Range("AB8").Value = '14:16'
